# Añadiendo nombre de modelo: new column F "Modelo" with the pipeline repr.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell (A1) onto the
# new header cell F1, then set its text.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "Modelo"

# New data cell for row 2 (plain, unstyled like the other row-2 cells).
$ws.Range("F2").Value = "Pipeline(steps=[('model', GradientBoostingRegressor(n_estimators=150))])"

$excel.CutCopyMode = 0
